# Week 17 data log: insert new player "R.Armstead" into the roster header
# row (after N.Cottrell, before D.Chark) on both the "Rushing" and
# "Receiving" sheets, pushing every subsequent header one column to the
# right and extending the "n" placeholder row out to the new last column.

$wb = $excel.ActiveWorkbook

$newHeaders = @(
    "T.Lawrence", "C.Beathard", "J.Robinson", "C.Hyde", "D.Ogunbowale",
    "T.Etienne", "N.Cottrell", "R.Armstead", "D.Chark", "M.Jones",
    "L.Shenault", "T.Johnson", "J.Agnew", "T.Austin", "T.Godwin",
    "L.Treadwell", "J.Mickens", "C.Manhertz", "J.O'Shaughnessy",
    "J.Hollister", "L.Farrell", "D.Arnold"
)

foreach ($sheetName in @("Rushing", "Receiving")) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Grow the new trailing column (W) so it inherits the bold/boxed
    # header style already used by the rest of row 1.
    $ws.Range("V1").Copy() | Out-Null
    $ws.Range("W1").PasteSpecial(-4122) | Out-Null

    # Write the header row, shifted to make room for the new player.
    for ($i = 0; $i -lt $newHeaders.Length; $i++) {
        $col = 2 + $i  # column B is index 2
        $ws.Cells.Item(1, $col).Value = $newHeaders[$i]
    }

    # Row 2 placeholder values ("n") now extend through the new column W.
    $ws.Range("B2:W2").Value = "n"
}

Write-Host "Week 17 data logged: R.Armstead added."
